$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F14").Value = 931
$wsExhibit.Range("F18").Value = 8547
$wsExhibit.Range("F24").Value = 237

# Sheet "本地生活" (Local Life)
$wsLocal = $wb.Worksheets.Item("本地生活")
$wsLocal.Range("F2").Value = 5589
$wsLocal.Range("F3").Value = 420

# Sheet "全部类型" (All Types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 5589
$wsAll.Range("F4").Value = 420
$wsAll.Range("F23").Value = 931
$wsAll.Range("F29").Value = 8547
$wsAll.Range("F39").Value = 237
